# Applies the "maj de la Part_List_RISM.xlsx" update:
# adds a new "Carte contacteur" parts-list section (rows 58-73) to Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Section title (bold, like the other "Carte ..." section headers) ---
$ws.Cells.Item(58, 2).Value = "Carte contacteur"
$ws.Cells.Item(58, 2).Font.Bold = $true

# --- Column header row (same layout as the other sections) ---
$ws.Cells.Item(59, 1).Value = "Type"
$ws.Cells.Item(59, 2).Value = "Valeur"
$ws.Cells.Item(59, 3).Value = "Référence EPSA"
$ws.Cells.Item(59, 4).Value = "Référence Farnell"
$ws.Cells.Item(59, 5).Value = "Emplacement"
$ws.Cells.Item(59, 6).Value = "Nb de pièces"
$ws.Cells.Item(59, 7).Value = "UdV"
$ws.Cells.Item(59, 8).Value = "Prix UdV"
$ws.Cells.Item(59, 9).Value = "Prix pour la carte"

# --- Data rows ---
$rows = @(
    @{ Row = 60; A = "Condensateur"; B = "10u";            E = "C1";        F = 1 },
    @{ Row = 61; A = "Condensateur"; B = "27p";            E = "C2,C4";     F = 2 },
    @{ Row = 62; A = "Condensateur"; B = "1u";             E = "C3";        F = 1 },
    @{ Row = 63; A = "Diode";        B = "DIODE ?R3";      E = "D1,D2,D3";  F = 3 },
    @{ Row = 64; A = "Relais";       B = "RELAIS-DPDT";    E = "K1,K2,K3";  F = 3 },
    @{ Row = 65; A = "Connecteur";   B = "KK2";             E = "P1";        F = 1 },
    @{ Row = 66; A = "Connecteur";   B = "KK4";             E = "P2";        F = 1 },
    @{ Row = 67; A = "Connecteur";   B = "Bornier 4";      E = "P3";        F = 1 },
    @{ Row = 68; A = "Transistor";   B = "FET_N";          E = "Q1,Q2,Q3";  F = 3 },
    @{ Row = 69; A = "Résistance";   B = "10k 0,25W";      E = "R1,R2";     F = 2 },
    @{ Row = 70; A = "Résistance";   B = "1k 0,25W";       E = "R3,R5,R7";  F = 3 },
    @{ Row = 71; A = "Résistance";   B = "220ohm 0,25W";   E = "R4,R6,R8";  F = 3 },
    @{ Row = 72; A = "LM";           B = "LM7805";         E = "U1";        F = 1 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# --- Row 73: "Logique" / U5, with B73 stored as a plain left-aligned number ---
$ws.Cells.Item(73, 1).Value = "Logique"
$ws.Cells.Item(73, 2).Value = 4001
$ws.Cells.Item(73, 2).HorizontalAlignment = -4131
$ws.Cells.Item(73, 5).Value = "U5"
$ws.Cells.Item(73, 6).Value = 1

# --- Scroll/selection state, as captured in the saved workbook ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F65").Select()
